$d = $word.ActiveDocument

# 0. Drop the "_GoBack" bookmark from the stray empty paragraph near the
#    end of the document; it is going to be re-created right after the
#    document title instead.
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBookmark = $d.Bookmarks.Item("_GoBack")
    $oldBookmark.Delete()
}

# 1. Locate the base title text "Regras de comunicação" (it is currently
#    followed by " (Communication rules)" inside the same paragraph).
$cCedilla = [char]0x00E7
$aTilde = [char]0x00E3
$baseText = "Regras de comunica" + $cCedilla + $aTilde + "o"

$findRange = $d.Content
$found = $findRange.Find.Execute($baseText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $bookmarkPos = $findRange.End

    # 2. Re-create the "_GoBack" bookmark, collapsed right after the base
    #    title text, while the trailing " (Communication rules)" text is
    #    still present (so this position is a plain in-text offset, not
    #    the end-of-paragraph / end-of-cell boundary, which this runtime
    #    refuses for a freshly-added collapsed bookmark).
    $bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)

    # 3. Find the paragraph that contains the title text, by scanning the
    #    document's Paragraphs collection (more reliable here than
    #    navigating ".Paragraphs" off an ad-hoc/Find range).
    $titlePara = $null
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $candidate = $d.Paragraphs.Item($i)
        if (($candidate.Range.Start -le $bookmarkPos) -and ($bookmarkPos -le $candidate.Range.End)) {
            $titlePara = $candidate
            break
        }
    }

    # 4. Remove the trailing " (Communication rules)" text (together with
    #    its italic runs and proofErr spell-check markers), leaving only
    #    "Regras de comunicação" followed by the bookmark.
    if ($titlePara -ne $null) {
        $cutEnd = $titlePara.Range.End - 1
        if ($cutEnd -gt $bookmarkPos) {
            $cutRange = $d.Range($bookmarkPos, $cutEnd)
            $cutRange.Delete()
        }
    }
}
